# semana 45 de 2025
# Adds week 44 (column AU) and week 45 (column AV) data to the weekly
# IRA/UCI revision sheet, and fixes a provider name in C52.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row: new week numbers 44 (AU1) and 45 (AV1) ---
# The existing week-number headers (D1..AT1) are stored as text, so force
# a text number format before assigning, otherwise "44"/"45" would be
# auto-converted to numbers.
$ws.Range("AU1").NumberFormat = "@"
$ws.Range("AV1").NumberFormat = "@"
$ws.Range("AU1").Value = "44"
$ws.Range("AV1").Value = "45"
$ws.Range("AU1").Font.Bold = $true
$ws.Range("AV1").Font.Bold = $true
$ws.Range("AU1").HorizontalAlignment = -4108
$ws.Range("AV1").HorizontalAlignment = -4108

# --- Row 41 also gains an AT41 (week 43) value that was previously blank ---
$ws.Range("AT41").Value = 0

# --- Week 44 (column AU) values per row ---
$au = @{
    2 = 0; 5 = 0; 6 = 1; 7 = 0; 8 = 0; 9 = 0; 10 = 0; 12 = 0; 14 = 0; 15 = 0;
    23 = 0; 24 = 0; 25 = 0; 28 = 1; 29 = 0; 30 = 1; 31 = 0; 34 = 0; 35 = 3;
    36 = 0; 37 = 0; 38 = 0; 40 = 0; 41 = 0; 42 = 0; 43 = 0; 44 = 0; 45 = 0;
    46 = 0; 47 = 0; 48 = 0; 49 = 0; 50 = 0; 51 = 0; 53 = 0; 54 = 0; 55 = 0;
    56 = 0; 57 = 0; 58 = 0
}

foreach ($row in $au.Keys) {
    $ws.Range("AU$row").Value = $au[$row]
}

# --- Week 45 (column AV) values per row ---
$av = @{
    2 = 0; 3 = 0; 5 = 0; 6 = 3; 7 = 0; 8 = 0; 9 = 0; 12 = 0; 13 = 0; 14 = 0;
    15 = 0; 16 = 0; 17 = 0; 22 = 0; 23 = 0; 25 = 0; 26 = 0; 28 = 0; 29 = 1;
    30 = 1; 31 = 0; 35 = 6; 36 = 0; 37 = 0; 38 = 0; 41 = 0; 42 = 0; 43 = 0;
    45 = 0; 46 = 0; 47 = 0; 48 = 0; 49 = 0; 50 = 0; 53 = 0; 54 = 0; 55 = 0;
    56 = 0; 57 = 0; 58 = 0
}

foreach ($row in $av.Keys) {
    $ws.Range("AV$row").Value = $av[$row]
}

# --- Update dimension-relevant / provider name fix ---
$ws.Range("C52").Value = "COOMEVA EXPERIENCIA MEDICA SAS"
